# Scheduled runner update: refresh market-price derived columns (H:N) on CRP and CUL leve tables,
# and clear stale cached price data (H:N) for ALC leve rows 125-141 (no longer tracked / out of range).

$wb = $excel.ActiveWorkbook

# --- ALC: clear stale cached columns H:N for rows 125-141 ---
$wsALC = $wb.Worksheets.Item("ALC")
$wsALC.Range("H125:N141").ClearContents()

# --- CRP & CUL: refresh computed price/profit columns with latest market data ---
$wsCRP = $wb.Worksheets.Item("CRP")
$wsCUL = $wb.Worksheets.Item("CUL")

$wsCRP.Range("H31").Value2 = 2489.2222
$wsCRP.Range("I31").Value2 = 1908.36
$wsCRP.Range("J31").Value2 = 9750
$wsCRP.Range("K31").Value2 = 1908.36
$wsCRP.Range("L31").Value2 = 9750
$wsCRP.Range("M31").Value2 = -1613.36
$wsCRP.Range("N31").Value2 = -10340
$wsCRP.Range("H34").Value2 = 2489.2222
$wsCRP.Range("I34").Value2 = 1908.36
$wsCRP.Range("J34").Value2 = 9750
$wsCRP.Range("K34").Value2 = 1908.36
$wsCRP.Range("L34").Value2 = 9750
$wsCRP.Range("M34").Value2 = -1706.36
$wsCRP.Range("N34").Value2 = -10154
$wsCRP.Range("H58").Value2 = 1068277
$wsCRP.Range("I58").Value2 = 2082.4333
$wsCRP.Range("J58").Value2 = 2949797
$wsCRP.Range("K58").Value2 = 2082.4333
$wsCRP.Range("L58").Value2 = 2949797
$wsCRP.Range("M58").Value2 = -1879.4333
$wsCRP.Range("N58").Value2 = -2950203
$wsCRP.Range("H122").Value2 = 52634330
$wsCRP.Range("I122").Value2 = 166667840
$wsCRP.Range("J122").Value2 = 3476.4614
$wsCRP.Range("K122").Value2 = 500003520
$wsCRP.Range("L122").Value2 = 10429.3842
$wsCRP.Range("M122").Value2 = -500001070
$wsCRP.Range("N122").Value2 = -15329.3842
$wsCRP.Range("H132").Value2 = 3155
$wsCRP.Range("I132").Value2 = 1965.3334
$wsCRP.Range("J132").Value2 = 3749.8333
$wsCRP.Range("K132").Value2 = 5896.0002
$wsCRP.Range("L132").Value2 = 11249.4999
$wsCRP.Range("M132").Value2 = -3366.0002
$wsCRP.Range("N132").Value2 = -16309.4999
$wsCRP.Range("H133").Value2 = 22323.125
$wsCRP.Range("J133").Value2 = 22323.125
$wsCRP.Range("L133").Value2 = 22323.125
$wsCRP.Range("N133").Value2 = -27383.125
$wsCRP.Range("H135").Value2 = 34946.668
$wsCRP.Range("J135").Value2 = 34946.668
$wsCRP.Range("L135").Value2 = 34946.668
$wsCRP.Range("N135").Value2 = -45086.668
$wsCRP.Range("H136").Value2 = 1068277
$wsCRP.Range("I136").Value2 = 2082.4333
$wsCRP.Range("J136").Value2 = 2949797
$wsCRP.Range("K136").Value2 = 6247.2999
$wsCRP.Range("L136").Value2 = 8849391
$wsCRP.Range("M136").Value2 = -3697.2999
$wsCRP.Range("N136").Value2 = -8854491
$wsCRP.Range("H137").Value2 = 25064.572
$wsCRP.Range("J137").Value2 = 25064.572
$wsCRP.Range("L137").Value2 = 25064.572
$wsCRP.Range("N137").Value2 = -35264.572
$wsCRP.Range("H138").Value2 = 40568
$wsCRP.Range("J138").Value2 = 40568
$wsCRP.Range("L138").Value2 = 40568
$wsCRP.Range("N138").Value2 = -50848
$wsCUL.Range("H5").Value2 = 1990.5385
$wsCUL.Range("J5").Value2 = 2686.5557
$wsCUL.Range("L5").Value2 = 8059.6671
$wsCUL.Range("N5").Value2 = -8283.667099999999
$wsCUL.Range("H68").Value2 = 1028.6842
$wsCUL.Range("I68").Value2 = 759
$wsCUL.Range("J68").Value2 = 1224.8182
$wsCUL.Range("K68").Value2 = 2277
$wsCUL.Range("L68").Value2 = 3674.4546
$wsCUL.Range("M68").Value2 = -1466
$wsCUL.Range("N68").Value2 = -5296.4546
$wsCUL.Range("H71").Value2 = 1028.6842
$wsCUL.Range("I71").Value2 = 759
$wsCUL.Range("J71").Value2 = 1224.8182
$wsCUL.Range("K71").Value2 = 6831
$wsCUL.Range("L71").Value2 = 11023.3638
$wsCUL.Range("M71").Value2 = -2775
$wsCUL.Range("N71").Value2 = -19135.3638
$wsCUL.Range("H122").Value2 = 3492.8
$wsCUL.Range("J122").Value2 = 3818.0813
$wsCUL.Range("L122").Value2 = 34362.7317
$wsCUL.Range("N122").Value2 = -39262.7317
$wsCUL.Range("H132").Value2 = 2473.6365
$wsCUL.Range("I132").Value2 = 1150
$wsCUL.Range("K132").Value2 = 10350
$wsCUL.Range("M132").Value2 = -7820
$wsCUL.Range("H135").Value2 = 1990.5385
$wsCUL.Range("J135").Value2 = 2686.5557
$wsCUL.Range("L135").Value2 = 24179.0013
$wsCUL.Range("N135").Value2 = -29249.0013
